# Update crypto price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.804.21"
$ws.Range("E2").Value = "  -2.85%  "
$ws.Range("D3").Value = "1.964.61"
$ws.Range("E3").Value = "  -2.14%  "
$ws.Range("D4").Value = "'1.015"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'323.59"
$ws.Range("E5").Value = "  -2.55%  "
$ws.Range("D6").Value = "'1.013"
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("D7").Value = "'0.4770"
$ws.Range("E7").Value = "  -5.08%  "
$ws.Range("D8").Value = "'0.4039"
$ws.Range("E8").Value = "  -5.42%  "
$ws.Range("D9").Value = "'53.68"
$ws.Range("E9").Value = "  -2.17%  "
$ws.Range("D10").Value = "'0.08508"
$ws.Range("E10").Value = "  -7.23%  "
$ws.Range("D11").Value = "'1.060"
$ws.Range("E11").Value = "  -6.06%  "
$ws.Range("D12").Value = "'22.42"
$ws.Range("E12").Value = "  -4.97%  "
$ws.Range("D13").Value = "1.991.11"
$ws.Range("E13").Value = "  -1.87%  "
$ws.Range("D14").Value = "'7.675"
$ws.Range("E14").Value = "  -5.80%  "
$ws.Range("D15").Value = "'6.249"
$ws.Range("E15").Value = "  -4.65%  "
$ws.Range("D17").Value = "'90.36"
$ws.Range("E17").Value = "  -5.41%  "
$ws.Range("D18").Value = "'0.00001067"
$ws.Range("E18").Value = "  -5.10%  "
$ws.Range("D19").Value = "'0.06616"
$ws.Range("E19").Value = "  -1.05%  "
$ws.Range("D20").Value = "'18.66"
$ws.Range("E20").Value = "  -6.31%  "
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("D22").Value = "'5.776"
$ws.Range("E22").Value = "  -3.60%  "
$ws.Range("D23").Value = "28.810.42"
$ws.Range("E23").Value = "  -2.84%  "
$ws.Range("D24").Value = "'11.55"
$ws.Range("E24").Value = "  -4.46%  "
$ws.Range("D25").Value = "'2.297"
$ws.Range("E25").Value = "  +0.66%  "
$ws.Range("D26").Value = "2.203.01"
$ws.Range("E26").Value = "  -2.53%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "'153.58"
$ws.Range("E27").Value = "  -3.66%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'20.22"
$ws.Range("E28").Value = "  -2.97%  "
$ws.Range("D29").Value = "'5.990"
$ws.Range("E29").Value = "  -6.81%  "
$ws.Range("D30").Value = "'2.162"
$ws.Range("E30").Value = "  -7.25%  "
$ws.Range("D31").Value = "'124.39"
$ws.Range("E31").Value = "  -3.55%  "
$ws.Range("D32").Value = "'1.005"
$ws.Range("E32").Value = "  -5.80%  "
$ws.Range("D33").Value = "'0.09642"
$ws.Range("E33").Value = "  -3.26%  "
$ws.Range("D34").Value = "'1.465"
$ws.Range("E34").Value = "  -7.70%  "
$ws.Range("D35").Value = "'5.680"
$ws.Range("E35").Value = "  -3.19%  "
$ws.Range("D36").Value = "'3.695"
$ws.Range("E36").Value = "  -3.08%  "
$ws.Range("D37").Value = "'0.02346"
$ws.Range("E37").Value = "  -5.51%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "'0.06245"
$ws.Range("E38").Value = "  -2.43%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "'1.271"
$ws.Range("E39").Value = "  -3.75%  "
$ws.Range("D40").Value = "'8.751"
$ws.Range("E40").Value = "  -8.53%  "
$ws.Range("D41").Value = "'0.6243"
$ws.Range("E41").Value = "  -5.50%  "
$ws.Range("D42").Value = "'11.17"
$ws.Range("E42").Value = "  -5.33%  "
$ws.Range("D43").Value = "'1.012"
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("E44").Value = "  -7.52%  "
$ws.Range("D45").Value = "'1.326"
$ws.Range("E45").Value = "  +2.33%  "
$ws.Range("D46").Value = "'0.5953"
$ws.Range("E46").Value = "  -6.82%  "
$ws.Range("D47").Value = "'12.90"
$ws.Range("E47").Value = "  -5.52%  "
$ws.Range("D48").Value = "'2.084"
$ws.Range("E48").Value = "  -6.09%  "
$ws.Range("D49").Value = "'3.430"
$ws.Range("E49").Value = "  -2.94%  "
$ws.Range("E50").Value = "  +0.62%  "
$ws.Range("D51").Value = "'2.104"
$ws.Range("E51").Value = "  +3.81%  "
